# Add 2022-Q3 data:
#  1) Insert a new "2022-Q3" row into the "总计" (total) summary sheet.
#  2) Insert a brand-new "2022-Q3" worksheet (holding-level detail) right
#     after "总计", pushing every later quarter sheet one slot down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: the row index column (A) never changes -- only the
#    date/count/value columns (B/C/D) cascade down by one quarter, with
#    the new 2022-Q3 figures landing in row 2 and row 8 newly created to
#    hold what used to be the last row (2020-Q4).
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Create row 8 as a copy of row 7 (same index-column style), then fill
# in the content cascade from the bottom up so nothing is clobbered
# before it's read.
$wsTotal.Range("A7:D7").Copy()
$wsTotal.Range("A8:D8").PasteSpecial(-4122)  # xlPasteFormats
$wsTotal.Cells.Item(8, 1).Value = 6

$cascade = @(
  @("2022-Q3", 14,   0.86),
  @("2022-Q2", 8,    0.57),
  @("2022-Q1", 13,   2.34),
  @("2021-Q4", 4,    0.11),
  @("2021-Q3", 2,    0.02),
  @("2021-Q2", 3,    0.05),
  @("2020-Q4", 2,    0.03)
)
for ($i = 0; $i -lt $cascade.Count; $i++) {
    $r = 2 + $i
    $row = $cascade[$i]
    $wsTotal.Cells.Item($r, 2).Value = $row[0]
    $wsTotal.Cells.Item($r, 3).Value = $row[1]
    $wsTotal.Cells.Item($r, 4).Value = $row[2]
}

# ---------------------------------------------------------------------
# 2) New "2022-Q3" worksheet, inserted right after "总计".
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Add($wsTotal.Next)
$wsQ3.Name = "2022-Q3"

# Header row (bold, bordered like the other quarter sheets).
$wsOldQ2 = $wb.Worksheets.Item("2022-Q2")
$wsOldQ2.Range("B1:H1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $wsQ3.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# Data rows.
$data = @(
  @("010744", "工银灵动价值混合A",       "11.46", "75.15", "2.66", "0.3048", 7),
  @("004350", "汇丰晋信价值先锋股票A",   "5.31",  "94.44", "3.91", "0.2076", 3),
  @("481008", "工银大盘蓝筹混合",        "4.05",  "66.88", "2.02", "0.0818", 10),
  @("481013", "工银消费服务混合A",       "2.76",  "70.46", "2.43", "0.0671", 9),
  @("007832", "博道伍佰智航股票C",       "6.02",  "88.19", "1.03", "0.0620", 5),
  @("007831", "博道伍佰智航股票A",       "3.13",  "88.19", "1.03", "0.0322", 5),
  @("920008", "中金进取回报灵活配置混合A", "1.59", "87.87", "1.96", "0.0312", 10),
  @("010745", "工银灵动价值混合C",       "0.87",  "75.15", "2.66", "0.0231", 7),
  @("920928", "中金进取回报灵活配置混合C", "1.13", "87.87", "1.96", "0.0221", 10),
  @("011179", "浙商智选食品饮料股票A",   "0.14",  "91.42", "8.46", "0.0118", 1),
  @("519987", "长信恒利优势混合",        "0.21",  "87.52", "4.60", "0.0097", 4),
  @("011180", "浙商智选食品饮料股票C",   "0.08",  "91.42", "8.46", "0.0068", 1),
  @("015364", "汇丰晋信价值先锋股票C",   "0.02",  "94.44", "3.91", "0.0008", 3),
  @("011475", "工银消费服务混合C",       "0.02",  "70.46", "2.43", "0.0005", 9)
)

$rowCount = $data.Count
$lastRow = 1 + $rowCount

# Index column A (numeric, bold+bordered style like the sibling sheets).
$wsOldQ2.Range("A2").Copy()
$wsQ3.Range("A2:A$lastRow").PasteSpecial(-4122)  # xlPasteFormats
for ($i = 0; $i -lt $rowCount; $i++) {
    $wsQ3.Cells.Item(2 + $i, 1).Value = $i
}

# Text columns B:G must stay text (fund codes with leading zeros, and
# numeric-looking figures stored as strings) -- force text format first.
$wsQ3.Range("B2:G$lastRow").NumberFormat = "@"
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = 2 + $i
    $row = $data[$i]
    $wsQ3.Cells.Item($r, 2).Value = $row[0]
    $wsQ3.Cells.Item($r, 3).Value = $row[1]
    $wsQ3.Cells.Item($r, 4).Value = $row[2]
    $wsQ3.Cells.Item($r, 5).Value = $row[3]
    $wsQ3.Cells.Item($r, 6).Value = $row[4]
    $wsQ3.Cells.Item($r, 7).Value = $row[5]
    $wsQ3.Cells.Item($r, 8).Value = $row[6]
}

Write-Host "2022-Q3 sheet and totals row added"
